$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("pull_down1","pull_down1_choice","pull_down2","pull_down2_choice","pull_down3","pull_down3_choice","input1","input2","input3")
$cols = @("I","J","K","L","M","N","O","P","Q")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $target = $ws.Range($cols[$i] + "1")
    $ws.Range("H1").Copy($target)
    $target.Value = $headers[$i]
}
